# Add two new "Manfred Man" rows (to 2B and to 3B) right after the header
# rows, pushing all the existing code rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at rows 3 and 4 (everything currently there,
# and below, shifts down by two rows).
$ws.Rows("3:4").Insert()

# Row 3: Manfred Man to 2B
$ws.Range("A3").Value = "Manfred Man to 2B"
$ws.Range("B3").Value = "B1"
$ws.Range("C3").Value = "M2"
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = "NoPitch"

# Row 4: Manfred Man to 3B
$ws.Range("A4").Value = "Manfred Man to 3B"
$ws.Range("B4").Value = "B1"
$ws.Range("C4").Value = "M3"
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = "NoPitch"

# Match the refreshed view state recorded in the workbook after the edit.
$ws.Range("E5").Select()
